$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.678.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.164.42"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.41"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.85"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.40%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.165.14"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.05%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.64"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.480"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.76%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.16"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.665.73"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.740.31"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.02%  "

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.114"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.163.69"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.07"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "484.28"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.90"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.82"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.42%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.25"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.78"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.31%  "

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.15"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.22%  "

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.122"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.75%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.95"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.50%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.16"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.86"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.25"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0748"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "460.34"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -8.46%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0405"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.02%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.87%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.898.91"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.48%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.07%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.96%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.69"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.09%  "

